$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.047.09"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.566.01"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("E10").Value = "  +1.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0859"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "1.558.16"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.521"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "27.044.87"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "216.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "154.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("E30").Value = "  +4.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("E32").Value = "  +4.38%  "
$ws.Range("D33").Value = "1.426.59"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("E34").Value = "  +12.39%  "
$ws.Range("E35").Value = "  +1.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.77%  "
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").Value = "1.702.69"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.52%  "
